# Apply "finish the others mapdata" edits to the Scene sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scene")

# B6: mine tunnel name corrected to "mine vein foothill"
$ws.Range("B6").Value = "矿脉山脚"

# TilePath column (I) filled in for rows that still had the placeholder "default"
$ws.Range("I6").Value = "orevalley"
$ws.Range("I12").Value = "gerdin"
$ws.Range("I19").Value = "fogvalley"
$ws.Range("I20").Value = "woodviliage"
$ws.Range("I21").Value = "riverside"
$ws.Range("I22").Value = "moonforest"

# Move the active selection to I22, matching the last edited cell
$ws.Range("I22").Select()
